# Applies the "Add files via upload" revision to parte1_pl6.xlsx:
#  - Column F (rows 2-21) now sums the residence-time columns C:E instead
#    of the utilisation columns G:I ("Tpo. Respuesta" = suma de tiempos,
#    not de utilizaciones).
#  - The newly-vacated utilisation numbers are no longer referenced, and a
#    stray formatted-but-empty column K (rows 2-21) appears (format carried
#    along with the F column edit).
#  - The small summary table (rows 44-64) gets a new "Tpo. Respuesta"
#    column C, built from the recomputed F values; the header cell that
#    used to read "TRES" becomes "Tpo. Respuesta", and the label/value pair
#    that used to live in C44/C45 ("TRES" / "SYSTEM_RESPONSE_TIME - Z") is
#    relocated out to O24/O25.
#  - Chart "Gr\u00e1fico 6" is moved from below the table to beside it.
#  - The view is scrolled/zoomed to where the author was last looking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-point the response-time column at the residence times ---------
# Mirrors the author's fill-down: F2 keeps its own literal formula, F3:F21
# carry on as one shared formula, all now =SUM(Cn:En) instead of =SUM(Gn:In).
$ws.Range("F2:F21").Formula = "=SUM(C2:E2)"

# --- 2. Stray number-formatted (but empty) K column alongside it ---------
$ws.Range("K2:K21").NumberFormat = "0.00E+00"

# --- 3. Relocate the old TRES / SYSTEM_RESPONSE_TIME - Z pair to O24:O25 -
$ws.Range("O24").Value = "TRES"
$ws.Range("O25").Value = "SYSTEM_RESPONSE_TIME - Z"

# --- 4. Re-purpose the table header C44 + add a literal "C" column -------
$ws.Range("C44").Value = "Tpo. Respuesta"

for ($r = 2; $r -le 21; $r++) {
    $destRow = 45 + ($r - 2)
    $v = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($destRow, 3).Value = $v
}

# --- 5. Move "Gr\u00e1fico 6" up beside the table -------------------------
$cos = $ws.ChartObjects()
$chart6 = $cos.Item(6)

$fromCell = $ws.Cells.Item(32, 14)
$toCell = $ws.Cells.Item(47, 20)
$fromColOff = 217714 / 12700
$fromRowOff = 95250 / 12700
$toColOff = 211315 / 12700
$toRowOff = 151 / 12700

$newLeft = $fromCell.Left() + $fromColOff
$newTop = $fromCell.Top() + $fromRowOff
$newRight = $toCell.Left() + $toColOff
$newBottom = $toCell.Top() + $toRowOff

$chart6.Left = $newLeft
$chart6.Top = $newTop
$chart6.Width = $newRight - $newLeft
$chart6.Height = $newBottom - $newTop

# --- 6. Restore the author's last view (scroll position, zoom, selection) -
$win = $excel.ActiveWindow
$win.ScrollRow = 64
$win.ScrollColumn = 2
$win.Zoom = 90
$ws.Range("G45").Select()
